# Auto-generated edit script: apply "Add data for 2022-11-11" updates
# across Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 110
$ws.Range("D3").Value = 126
$ws.Range("B9").Value = 352
$ws.Range("C9").Value = 452
$ws.Range("D9").Value = 383
$ws.Range("E9").Value = 426
$ws.Range("H9").Value = 420
$ws.Range("I9").Value = 471
$ws.Range("B10").Value = 1237
$ws.Range("C10").Value = 1456
$ws.Range("D10").Value = 1677
$ws.Range("E10").Value = 1982
$ws.Range("F10").Value = 1970
$ws.Range("G10").Value = 855
$ws.Range("I10").Value = 796
$ws.Range("B11").Value = 1714
$ws.Range("C11").Value = 2061
$ws.Range("D11").Value = 2282
$ws.Range("E11").Value = 2623
$ws.Range("F11").Value = 2666
$ws.Range("G11").Value = 1495
$ws.Range("H11").Value = 1231
$ws.Range("I11").Value = 1587

# --- Sheet 12: Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("C7").Value = 33
$ws.Range("I7").Value = 30
$ws.Range("C8").Value = 56
$ws.Range("E8").Value = 76
$ws.Range("C9").Value = 94
$ws.Range("E9").Value = 135
$ws.Range("I9").Value = 92

# --- Sheet 14: Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("D3").Value = 8
$ws.Range("D9").Value = 74

# --- Sheet 15: Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E7").Value = 61
$ws.Range("B8").Value = 193
$ws.Range("C8").Value = 293
$ws.Range("D8").Value = 485
$ws.Range("E8").Value = 581
$ws.Range("G8").Value = 156
$ws.Range("B9").Value = 239
$ws.Range("C9").Value = 345
$ws.Range("D9").Value = 552
$ws.Range("E9").Value = 658
$ws.Range("G9").Value = 236

# --- Sheet 16: Armour Square ---
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("D6").Value = 9
$ws.Range("D7").Value = 16

# --- Sheet 17: Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("F6").Value = 77
$ws.Range("F7").Value = 100

# --- Sheet 18: Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("E5").Value = 18
$ws.Range("E7").Value = 45

# --- Sheet 19: North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("D6").Value = 20
$ws.Range("D7").Value = 41
$ws.Range("E7").Value = 33
$ws.Range("D8").Value = 62
$ws.Range("E8").Value = 44

# --- Sheet 2: By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("D5").Value = 16
$ws.Range("F8").Value = 133
$ws.Range("G10").Value = 12
$ws.Range("D18").Value = 6
$ws.Range("B27").Value = 19
$ws.Range("E28").Value = 85
$ws.Range("F28").Value = 113
$ws.Range("C32").Value = 94
$ws.Range("E32").Value = 135
$ws.Range("I32").Value = 92
$ws.Range("D36").Value = 74
$ws.Range("I41").Value = 15
$ws.Range("I42").Value = 10
$ws.Range("F43").Value = 16
$ws.Range("E50").Value = 45
$ws.Range("B53").Value = 239
$ws.Range("C53").Value = 345
$ws.Range("D53").Value = 552
$ws.Range("E53").Value = 658
$ws.Range("G53").Value = 236
$ws.Range("H61").Value = 22
$ws.Range("I61").Value = 15
$ws.Range("F63").Value = 12
$ws.Range("D65").Value = 62
$ws.Range("E65").Value = 44
$ws.Range("F70").Value = 100
$ws.Range("I72").Value = 14
$ws.Range("D74").Value = 92
$ws.Range("E76").Value = 91
$ws.Range("F78").Value = 50
$ws.Range("E92").Value = 30
$ws.Range("I92").Value = 34
$ws.Range("D96").Value = 31
$ws.Range("B97").Value = 28
$ws.Range("H98").Value = 7
$ws.Range("B99").Value = 1714
$ws.Range("C99").Value = 2061
$ws.Range("D99").Value = 2282
$ws.Range("E99").Value = 2623
$ws.Range("F99").Value = 2666
$ws.Range("G99").Value = 1495
$ws.Range("H99").Value = 1231
$ws.Range("I99").Value = 1587

# --- Sheet 22: Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 1
$ws.Range("I5").Value = 6
$ws.Range("I6").Value = 15

# --- Sheet 25: Rush & Division ---
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("F5").Value = 40
$ws.Range("F6").Value = 50

# --- Sheet 26: Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("E8").Value = 55
$ws.Range("F8").Value = 62
$ws.Range("E9").Value = 85
$ws.Range("F9").Value = 113

# --- Sheet 3: Rogers Park ---
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("E7").Value = 73
$ws.Range("E8").Value = 91

# --- Sheet 31: River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("D6").Value = 81
$ws.Range("D7").Value = 92

# --- Sheet 33: West Loop ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I2").Value = 3
$ws.Range("I7").Value = 6
$ws.Range("E8").Value = 25
$ws.Range("E9").Value = 30
$ws.Range("I9").Value = 34

# --- Sheet 37: Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("B5").Value = 10
$ws.Range("B7").Value = 28

# --- Sheet 4: Edgewater ---
$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("B4").Value = 3
$ws.Range("B6").Value = 19

# --- Sheet 45: Calumet Heights ---
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 6

# --- Sheet 56: Wicker Park ---
$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("D6").Value = 28
$ws.Range("D7").Value = 31

# --- Sheet 58: New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 12

# --- Sheet 62: Printers Row ---
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("I4").Value = 9
$ws.Range("I6").Value = 14

# --- Sheet 65: Wrigleyville ---
$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("H5").Value = 1
$ws.Range("H7").Value = 7

# --- Sheet 66: Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 10

# --- Sheet 69: Avondale ---
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("G6").Value = 8
$ws.Range("G7").Value = 12

# --- Sheet 8: Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("F7").Value = 91
$ws.Range("F8").Value = 133

# --- Sheet 82: Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 16
